{"js": "// Change \"\u01af\u1edbc l\u01b0\u1ee3ng s\u1ed1 testcase: 100 testcase\" -> \"... : 90 testcase\"\n// (diff splits the old \"100 testcase\" run into \"90\" + \" testcase\").\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"testcase\") !== -1\n);\nif (!target) {\n  throw new Error('Could not find the \"testcase\" paragraph.');\n}\n\n// Replace just the number - this keeps the edited run's boundary separate\n// from the preceding \": \" run (unlike replacing the whole \"100 testcase\").\nconst numberHits = target.search(\"100\", { matchCase: true });\nnumberHits.load(\"items\");\nawait context.sync();\nif (numberHits.items.length === 0) {\n  throw new Error('Could not find \"100\" in the target paragraph.');\n}\nnumberHits.items[0].insertText(\"90\", \"Replace\");\nawait context.sync();\n\n// Split \"90 testcase\" into \"90\" + \" testcase\" runs, matching the diff,\n// by toggling italic off/on over just the trailing \" testcase\" text -\n// Word folds same-formatting runs back together on a plain text edit, so\n// this forces the run boundary the diff shows without changing the look.\nconst suffixHits = target.search(\" testcase\", { matchCase: true });\nsuffixHits.load(\"items\");\nawait context.sync();\nconst suffix = suffixHits.items[suffixHits.items.length - 1];\nsuffix.font.italic = false;\nawait context.sync();\nsuffix.font.italic = true;\nawait context.sync();\n", "ps1": "# Change \"\u01af\u1edbc l\u01b0\u1ee3ng s\u1ed1 testcase: 100 testcase\" -> \"... : 90 testcase\"\n# (diff splits the old \"100 testcase\" run into \"90\" + \" testcase\").\n\n$d = $word.ActiveDocument\n\n# Locate the run's text (\"100 testcase\") without mutating the document yet -\n# Find.Execute on a throw-away range just repositions that range to the hit.\n$hit = $d.Content\n$found = $hit.Find.Execute(\"100 testcase\")\nif (-not $found) {\n    throw \"Could not find '100 testcase' in the document.\"\n}\n\n# Narrow down to just the number \"100\" at the start of that run.\n$numStart = $hit.Start\n$numEnd = $numStart + 3\n$numRange = $d.Range($numStart, $numEnd)\n\n# Toggling Italic off/on around the Text assignment keeps this run split\n# from its neighbours (otherwise assigning Range.Text here folds the whole\n# paragraph's like-formatted runs back into a single run).\n$numRange.Font.Italic = $false\n$numRange.Text = \"90\"\n$numRange.Font.Italic = $true\n"}
